$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume table with the latest scraped values.
# Numeric-looking text values (prices, percentages) are written with a
# leading apostrophe so Excel keeps them as text (matching the source data,
# which stores them as inline strings, not numbers), then the style is reset
# back to Normal so no stray number-format / quote-prefix formatting sticks.

$ws.Range("D2").Value = "'262.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.71%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.03%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.699"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.30%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06087"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.14%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.700"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.47%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'-0.03%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9083"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.78%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.05194"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'9.29%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07091"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.05%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.03132"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.01%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09039"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.17%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001527"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.74%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006154"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.08%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005993"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.07%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.452"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.09%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.165"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.28%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.167"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.56%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3072"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.21%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D20").Value = "'0.1406"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.03%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.29%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.086"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.02%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04246"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.49%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-3.26%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004056"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.73%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E27").Value = "'23.03%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03939"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.72%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.08%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.004171"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.24%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002109"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.01145"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-29.86%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005079"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.66%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'0.2510"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'51.72%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("E50").Style = "Normal"
